# Translate the hard-coded English subtitle lines to Italian.
# Each old string is unique in the document, so Find locates the
# single matching run; we then set Range.Text directly (rather than
# passing a Replacement string to Find.Execute) so Word's
# autocorrect/smart-quote substitution doesn't mangle the straight
# apostrophes used in the source text.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        $rng.Text = $newText
        Write-Output "Replaced: $oldText"
    } else {
        Write-Output "NOT FOUND: $oldText"
    }
}

Replace-ExactText "** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino" "** il dialogo inizia a 47 secondi, quindi ne ho aggiunti 28 in ogni momento dove necessari. -John Argentino"

Replace-ExactText "Mike, astonished, asks: 'How old are they?' Fil," "Mike, stupito, chiede 'Quanti anni hanno?' Fil,"

Replace-ExactText "being a playful mathematician, answers" "essendo un matematico giocoso, risponde"

Replace-ExactText "'You tell me! I'll give you a hint: if you" "'Dimmelo tu! Ti darà un indizio: se tu"

Replace-ExactText "multiply the three ages together you" "moltiplichi le tre età"

Replace-ExactText "get 36.' Mike takes sometimes to think" "ottieni 36.' Mike impiega un po' di tempo a pensare"

Replace-ExactText "and says: 'I'm sorry Fil, but I do need" "e dice: 'Mi spiace Fil, ma mi serve"

Replace-ExactText "another hint. So Fil tells Mike:" "un altro suggerimento. Quindi Fil dice a Mike:"
